# Auto-generated edits applying the diff to Seraph_Profits workbook
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(32, 8).Value2 = 39999
$ws.Cells.Item(32, 9).Value2 = 30000
$ws.Cells.Item(32, 10).Value2 = 49998
$ws.Cells.Item(32, 11).Value2 = 30000
$ws.Cells.Item(32, 12).Value2 = 49998
$ws.Cells.Item(32, 13).Value2 = -29674
$ws.Cells.Item(32, 14).Value2 = -50650
$ws.Cells.Item(64, 8).Value2 = 3174.25
$ws.Cells.Item(64, 9).Value2 = 3174.25
$ws.Cells.Item(64, 11).Value2 = 3174.25
$ws.Cells.Item(64, 13).Value2 = -2926.25
$ws.Cells.Item(67, 8).Value2 = 3174.25
$ws.Cells.Item(67, 9).Value2 = 3174.25
$ws.Cells.Item(67, 11).Value2 = 3174.25
$ws.Cells.Item(67, 13).Value2 = -2316.25
$ws.Cells.Item(97, 8).Value2 = 3216.5
$ws.Cells.Item(97, 10).Value2 = 3216.5
$ws.Cells.Item(97, 12).Value2 = 9649.5
$ws.Cells.Item(97, 14).Value2 = -10641.5
$ws.Cells.Item(99, 8).Value2 = 581.5
$ws.Cells.Item(99, 9).Value2 = 164
$ws.Cells.Item(99, 10).Value2 = 999
$ws.Cells.Item(99, 11).Value2 = 492
$ws.Cells.Item(99, 12).Value2 = 2997
$ws.Cells.Item(99, 13).Value2 = 1006
$ws.Cells.Item(99, 14).Value2 = -5993
$ws.Cells.Item(103, 8).Value2 = 3500
$ws.Cells.Item(103, 9).Value2 = 4000
$ws.Cells.Item(103, 11).Value2 = 12000
$ws.Cells.Item(103, 13).Value2 = -11414
$ws.Cells.Item(107, 8).Value2 = 37038116
$ws.Cells.Item(107, 9).Value2 = 43479390
$ws.Cells.Item(107, 10).Value2 = 787.5
$ws.Cells.Item(107, 11).Value2 = 43479390
$ws.Cells.Item(107, 12).Value2 = 787.5
$ws.Cells.Item(107, 13).Value2 = -43477470
$ws.Cells.Item(107, 14).Value2 = -4627.5
$ws.Cells.Item(109, 8).Value2 = 40000
$ws.Cells.Item(109, 10).Value2 = 40000
$ws.Cells.Item(109, 12).Value2 = 40000
$ws.Cells.Item(109, 14).Value2 = -42774
$ws.Cells.Item(132, 8).Value2 = 2290.077
$ws.Cells.Item(132, 9).Value2 = 1299.3043
$ws.Cells.Item(132, 11).Value2 = 3897.9129
$ws.Cells.Item(132, 13).Value2 = -1367.9129
$ws.Cells.Item(135, 8).Value2 = 1525.9333
$ws.Cells.Item(135, 9).Value2 = 1303.2222
$ws.Cells.Item(135, 11).Value2 = 11728.9998
$ws.Cells.Item(135, 13).Value2 = -9193.9998

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(5, 8).Value2 = 630.625
$ws.Cells.Item(5, 9).Value2 = 499.16666
$ws.Cells.Item(5, 10).Value2 = 1025
$ws.Cells.Item(5, 11).Value2 = 499.16666
$ws.Cells.Item(5, 12).Value2 = 1025
$ws.Cells.Item(5, 13).Value2 = -387.16666
$ws.Cells.Item(5, 14).Value2 = -1249
$ws.Cells.Item(8, 8).Value2 = 51500000
$ws.Cells.Item(8, 9).Value2 = 51500000
$ws.Cells.Item(8, 11).Value2 = 51500000
$ws.Cells.Item(8, 13).Value2 = -51499856
$ws.Cells.Item(61, 8).Value2 = 8268.444
$ws.Cells.Item(61, 9).Value2 = 8677
$ws.Cells.Item(61, 10).Value2 = 5000
$ws.Cells.Item(61, 11).Value2 = 8677
$ws.Cells.Item(61, 12).Value2 = 5000
$ws.Cells.Item(61, 13).Value2 = -8465
$ws.Cells.Item(61, 14).Value2 = -5424
$ws.Cells.Item(74, 8).Value2 = 2984.5715
$ws.Cells.Item(74, 9).Value2 = 2225
$ws.Cells.Item(74, 11).Value2 = 2225
$ws.Cells.Item(74, 13).Value2 = -1351
$ws.Cells.Item(77, 8).Value2 = 2984.5715
$ws.Cells.Item(77, 9).Value2 = 2225
$ws.Cells.Item(77, 11).Value2 = 11125
$ws.Cells.Item(77, 13).Value2 = -6757
$ws.Cells.Item(97, 8).Value2 = 670.1429000000001
$ws.Cells.Item(97, 9).Value2 = 562.75
$ws.Cells.Item(97, 11).Value2 = 562.75
$ws.Cells.Item(97, 13).Value2 = -66.75
$ws.Cells.Item(136, 8).Value2 = 8268.444
$ws.Cells.Item(136, 9).Value2 = 8677
$ws.Cells.Item(136, 10).Value2 = 5000
$ws.Cells.Item(136, 11).Value2 = 26031
$ws.Cells.Item(136, 12).Value2 = 15000
$ws.Cells.Item(136, 13).Value2 = -23481
$ws.Cells.Item(136, 14).Value2 = -20100

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(4, 8).Value2 = 630.625
$ws.Cells.Item(4, 9).Value2 = 499.16666
$ws.Cells.Item(4, 10).Value2 = 1025
$ws.Cells.Item(4, 11).Value2 = 499.16666
$ws.Cells.Item(4, 12).Value2 = 1025
$ws.Cells.Item(4, 13).Value2 = -384.16666
$ws.Cells.Item(4, 14).Value2 = -1255
$ws.Cells.Item(26, 8).Value2 = 20332.666
$ws.Cells.Item(26, 9).Value2 = 20332.666
$ws.Cells.Item(26, 11).Value2 = 20332.666
$ws.Cells.Item(26, 13).Value2 = -20040.666
$ws.Cells.Item(29, 8).Value2 = 1899.5
$ws.Cells.Item(29, 9).Value2 = 800
$ws.Cells.Item(29, 10).Value2 = 2999
$ws.Cells.Item(29, 11).Value2 = 800
$ws.Cells.Item(29, 12).Value2 = 2999
$ws.Cells.Item(29, 13).Value2 = -511
$ws.Cells.Item(29, 14).Value2 = -3577
$ws.Cells.Item(36, 8).Value2 = 8019.5
$ws.Cells.Item(36, 9).Value2 = 998.5
$ws.Cells.Item(36, 11).Value2 = 998.5
$ws.Cells.Item(36, 13).Value2 = -464.5
$ws.Cells.Item(105, 8).Value2 = 3334830.5
$ws.Cells.Item(105, 9).Value2 = 4630642
$ws.Cells.Item(105, 10).Value2 = 2744
$ws.Cells.Item(105, 11).Value2 = 4630642
$ws.Cells.Item(105, 12).Value2 = 2744
$ws.Cells.Item(105, 13).Value2 = -4628895
$ws.Cells.Item(105, 14).Value2 = -6238
$ws.Cells.Item(134, 8).Value2 = 3035
$ws.Cells.Item(134, 9).Value2 = 2359.5454
$ws.Cells.Item(134, 11).Value2 = 7078.6362
$ws.Cells.Item(134, 13).Value2 = -4543.6362

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value2 = 62503170
$ws.Cells.Item(16, 9).Value2 = 62503170
$ws.Cells.Item(16, 11).Value2 = 62503170
$ws.Cells.Item(16, 13).Value2 = -62502883
$ws.Cells.Item(22, 8).Value2 = 64369
$ws.Cells.Item(22, 9).Value2 = 80549.87
$ws.Cells.Item(22, 11).Value2 = 80549.87
$ws.Cells.Item(22, 13).Value2 = -80199.87
$ws.Cells.Item(31, 8).Value2 = 2742
$ws.Cells.Item(31, 9).Value2 = 1665.85
$ws.Cells.Item(31, 11).Value2 = 1665.85
$ws.Cells.Item(31, 13).Value2 = -1370.85
$ws.Cells.Item(34, 8).Value2 = 2742
$ws.Cells.Item(34, 9).Value2 = 1665.85
$ws.Cells.Item(34, 11).Value2 = 1665.85
$ws.Cells.Item(34, 13).Value2 = -1463.85
$ws.Cells.Item(94, 8).Value2 = 2518.6
$ws.Cells.Item(94, 9).Value2 = 1949.5
$ws.Cells.Item(94, 10).Value2 = 2660.875
$ws.Cells.Item(94, 11).Value2 = 1949.5
$ws.Cells.Item(94, 12).Value2 = 2660.875
$ws.Cells.Item(94, 13).Value2 = -1498.5
$ws.Cells.Item(94, 14).Value2 = -3562.875
$ws.Cells.Item(113, 8).Value2 = 62503170
$ws.Cells.Item(113, 9).Value2 = 62503170
$ws.Cells.Item(113, 11).Value2 = 62503170
$ws.Cells.Item(113, 13).Value2 = -62501000

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(11, 8).Value2 = 758.5
$ws.Cells.Item(11, 9).Value2 = 513
$ws.Cells.Item(11, 10).Value2 = 1249.5
$ws.Cells.Item(11, 11).Value2 = 1539
$ws.Cells.Item(11, 12).Value2 = 3748.5
$ws.Cells.Item(11, 13).Value2 = -1399
$ws.Cells.Item(11, 14).Value2 = -4028.5
$ws.Cells.Item(14, 8).Value2 = 723.63635
$ws.Cells.Item(14, 9).Value2 = 723.63635
$ws.Cells.Item(14, 11).Value2 = 2170.90905
$ws.Cells.Item(14, 13).Value2 = -1997.90905
$ws.Cells.Item(106, 8).Value2 = 14000
$ws.Cells.Item(106, 10).Value2 = 14000
$ws.Cells.Item(106, 12).Value2 = 42000
$ws.Cells.Item(106, 14).Value2 = -43892
$ws.Cells.Item(122, 8).Value2 = 267.57144
$ws.Cells.Item(122, 9).Value2 = 260.25
$ws.Cells.Item(122, 11).Value2 = 2342.25
$ws.Cells.Item(122, 13).Value2 = 107.75
$ws.Cells.Item(131, 8).Value2 = 1406.238
$ws.Cells.Item(132, 8).Value2 = 3385.25
$ws.Cells.Item(132, 9).Value2 = 3091
$ws.Cells.Item(132, 11).Value2 = 27819
$ws.Cells.Item(132, 13).Value2 = -25289
$ws.Cells.Item(140, 8).Value2 = 1119.421
$ws.Cells.Item(140, 9).Value2 = 903.94446
$ws.Cells.Item(140, 11).Value2 = 2711.83338
$ws.Cells.Item(140, 13).Value2 = 2468.16662

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(41, 8).Value2 = 7500
$ws.Cells.Item(41, 9).Value2 = 7500
$ws.Cells.Item(41, 11).Value2 = 7500
$ws.Cells.Item(41, 13).Value2 = -7145
$ws.Cells.Item(122, 8).Value2 = 70495
$ws.Cells.Item(122, 9).Value2 = 3098.8
$ws.Cells.Item(122, 10).Value2 = 205287.4
$ws.Cells.Item(122, 11).Value2 = 9296.400000000001
$ws.Cells.Item(122, 12).Value2 = 615862.2
$ws.Cells.Item(122, 13).Value2 = -6846.400000000001
$ws.Cells.Item(122, 14).Value2 = -620762.2
$ws.Cells.Item(132, 8).Value2 = 3959
$ws.Cells.Item(132, 9).Value2 = 4083.762
$ws.Cells.Item(132, 11).Value2 = 12251.286
$ws.Cells.Item(132, 13).Value2 = -9721.286

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value2 = 7277.92
$ws.Cells.Item(22, 9).Value2 = 3095.6667
$ws.Cells.Item(22, 10).Value2 = 11138.462
$ws.Cells.Item(22, 11).Value2 = 3095.6667
$ws.Cells.Item(22, 12).Value2 = 11138.462
$ws.Cells.Item(22, 13).Value2 = -2800.6667
$ws.Cells.Item(22, 14).Value2 = -11728.462
$ws.Cells.Item(27, 8).Value2 = 7277.92
$ws.Cells.Item(27, 9).Value2 = 3095.6667
$ws.Cells.Item(27, 10).Value2 = 11138.462
$ws.Cells.Item(27, 11).Value2 = 3095.6667
$ws.Cells.Item(27, 12).Value2 = 11138.462
$ws.Cells.Item(27, 13).Value2 = -2988.6667
$ws.Cells.Item(27, 14).Value2 = -11352.462
$ws.Cells.Item(100, 8).Value2 = 7120.5
$ws.Cells.Item(100, 9).Value2 = 6995
$ws.Cells.Item(100, 11).Value2 = 6995
$ws.Cells.Item(100, 13).Value2 = -6454
$ws.Cells.Item(122, 8).Value2 = 5489.5884
$ws.Cells.Item(122, 9).Value2 = 3791.75
$ws.Cells.Item(122, 11).Value2 = 11375.25
$ws.Cells.Item(122, 13).Value2 = -8925.25
$ws.Cells.Item(136, 9).Value2 = 7751.5
$ws.Cells.Item(136, 10).Value2 = 8221.333000000001
$ws.Cells.Item(136, 11).Value2 = 23254.5
$ws.Cells.Item(136, 12).Value2 = 24663.999
$ws.Cells.Item(136, 13).Value2 = -20704.5
$ws.Cells.Item(136, 14).Value2 = -29763.999

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(54, 8).Value2 = 27487.5
$ws.Cells.Item(54, 10).Value2 = 29983.334
$ws.Cells.Item(54, 12).Value2 = 29983.334
$ws.Cells.Item(54, 14).Value2 = -31023.334
$ws.Cells.Item(56, 8).Value2 = 55000
$ws.Cells.Item(56, 9).Value2 = 55000
$ws.Cells.Item(56, 11).Value2 = 55000
$ws.Cells.Item(56, 13).Value2 = -54286
$ws.Cells.Item(108, 8).Value2 = 0
$ws.Cells.Item(108, 10).Value2 = 0
$ws.Cells.Item(108, 12).Value2 = 0
$ws.Cells.Item(108, 14).ClearContents()
$ws.Cells.Item(132, 8).Value2 = 1585.875
$ws.Cells.Item(132, 9).Value2 = 1336.2858
$ws.Cells.Item(132, 11).Value2 = 4008.8574
$ws.Cells.Item(132, 13).Value2 = -1478.8574
